# ---------------------------------------------------------------------------
# Update with Correct Forecast output
#
# - Renames Sheet1 to "Sales vs PO" and inserts an "Order Week" column
#   (old "ds" values), shifting "ds" forward one week and zeroing the
#   PO_Requested_Qty column (moved into the new "Weekly Growth" sheet).
# - Adds three new report sheets: "Weekly Growth", "Volume Insights",
#   and "Prediction Info" with the forecast summary output.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Sales vs PO" ------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column before the old "PO_Requested_Qty" column (C), pushing
# it to D. The new column C becomes "Order Week".
$ws1.Columns("C:C").Insert()
$ws1.Range("C1").Value = "Order Week"

$oldDs = @(45565, 45572, 45579, 45586, 45593, 45600, 45607, 45614, 45621, 45628, 45635, 45642, 45649)
$newDs = @(45571, 45578, 45585, 45592, 45599, 45606, 45613, 45620, 45627, 45634, 45641, 45648, 45655)

for ($i = 0; $i -lt $oldDs.Length; $i++) {
    $row = $i + 2
    $ws1.Range("A$row").Value = $newDs[$i]
    $ws1.Cells.Item($row, 3).Value = $oldDs[$i]
    $ws1.Cells.Item($row, 4).Value = 0
}

# New "Order Week" column should carry the same date formatting as "ds".
$ws1.Range("A2:A14").Copy()
$ws1.Range("C2:C14").PasteSpecial(-4122)

# --- Sheet 2: "Weekly Growth" ----------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$ws2.Range("A2").Value = 45572
$ws2.Range("B2").Value = 160
$ws2.Range("C2").Value = 0

$ws2.Range("A3").Value = 45586
$ws2.Range("B3").Value = 16
$ws2.Range("C3").Value = -90

# Match header / date styling to the "Sales vs PO" sheet.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws1.Range("A2:A3").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

# --- Sheet 3: "Volume Insights" ---------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 176
$ws3.Range("B2").Value = 88
$ws3.Range("C2").Value = 160
$ws3.Range("D2").Value = 16

$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# --- Sheet 4: "Prediction Info" ---------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 0

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws1.Range("A1").Select() | Out-Null
